# Insert a new weekly price record at row 26 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 26-54 down to 27-55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 26, shifting rows 26..54 -> 27..55.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44897
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = 300000000
$ws.Cells.Item(26, 7).Value = "Espárragos"
$ws.Cells.Item(26, 8).Value = "Verde"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 1100
$ws.Cells.Item(26, 11).Value = 1600
$ws.Cells.Item(26, 12).Value = 1600
$ws.Cells.Item(26, 13).Value = 1600
$ws.Cells.Item(26, 14).Value = "$/kilo"
$ws.Cells.Item(26, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(26, 16).Value = 1600
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
